$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column BB, mirroring column BA (header date row + data rows 1-82) ---
# First copy formatting only (so header cell BB1 gets the same bold/border/date-format
# style as BA1, and data cells BB2:BB82 stay unstyled, matching BA2:BA82).
$ws.Range("BA1:BA82").Copy()
$ws.Range("BB1:BB82").PasteSpecial(-4122)

# Then copy the values themselves (separately, so no new number-format style gets
# synthesised the way a combined "paste all" would).
$ws.Range("BA1:BA82").Copy()
$ws.Range("BB1:BB82").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# New header date for column BB (next forecast vintage, serial date 45986)
$ws.Range("BB1").Value = 45986

# Rows 72 onward take a revised forecast value, different from column BA
$ws.Range("BB72").Value = -0.5
$ws.Range("BB73").Value = -0.3
$ws.Range("BB74").Value = -0.25
$ws.Range("BB75").Value = -0.25
$ws.Range("BB76").Value = -0.25
$ws.Range("BB77").Value = -0.25
$ws.Range("BB78").Value = -0.25
$ws.Range("BB79").Value = -0.25
$ws.Range("BB80").Value = -0.25
$ws.Range("BB81").Value = -0.25
$ws.Range("BB82").Value = -0.25

# --- Append a new trailing row (83) for the extra forecast quarter ---
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = -0.25
